# Apply the changes described by the diff:
#  - Append a new data row (row 5) to the "speciesevents" sheet with a new
#    Station_ID "CTD" (adds a new shared string), utm_y/utm_x coordinates,
#    and Setup_date/Retrieval_date values formatted as dates.
#  - Update the active cell selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of camera-trap station data
$ws.Range("A5").Value = "CTD"
$ws.Range("B5").Value = 9134567
$ws.Range("C5").Value = 213394

# Apply the existing date number format before assigning the date value so
# the date cells reuse the workbook's existing style instead of creating a
# new, separate number format / cell style.
$ws.Range("D5:E5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("D5").Value = [DateTime]"2015-09-30"
$ws.Range("E5").Value = [DateTime]"2015-10-01"

# Update the selected cell shown when the workbook is reopened
$ws.Range("E10").Select()
